$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.20631253961366
$ws.Range("C2").Value = 6.549755449037868
$ws.Range("D2").Value = 6.005083849300747
$ws.Range("E2").Value = 11.1178129689675
$ws.Range("G2").Value = 3.668445937193999
$ws.Range("I2").Value = 25.86373005740542
$ws.Range("K2").Value = 10.67884887355966
$ws.Range("L2").Value = 9.933376744691756
$ws.Range("N2").Value = 19.97476744142771
$ws.Range("O2").Value = 26.58240624963684

$ws.Range("B3").Value = 12.97195362379746
$ws.Range("C3").Value = 6.475267126448523
$ws.Range("D3").Value = 5.889741399151301
$ws.Range("E3").Value = 11.12573610350912
$ws.Range("G3").Value = 3.670420806046276
$ws.Range("I3").Value = 25.93767183367571
$ws.Range("K3").Value = 10.51668087114148
$ws.Range("L3").Value = 9.920325693011829
$ws.Range("N3").Value = 20.03631681232589
$ws.Range("O3").Value = 26.64037423884268

$ws.Range("B4").Value = 12.82861737789182
$ws.Range("C4").Value = 6.42835473892077
$ws.Range("D4").Value = 5.819527884897889
$ws.Range("E4").Value = 11.13250006249051
$ws.Range("G4").Value = 3.671698395628471
$ws.Range("I4").Value = 25.98765482612441
$ws.Range("K4").Value = 10.41788693103296
$ws.Range("L4").Value = 9.914064216460321
$ws.Range("N4").Value = 20.07586443551667
$ws.Range("O4").Value = 26.68122505884

$ws.Range("B5").Value = 12.77042771902384
$ws.Range("C5").Value = 6.408950247951513
$ws.Range("D5").Value = 5.79111182881732
$ws.Range("E5").Value = 11.13573439053635
$ws.Range("G5").Value = 3.672235422906563
$ws.Range("I5").Value = 26.00917430303366
$ws.Range("K5").Value = 10.37787360434237
$ws.Range("L5").Value = 9.911955426378887
$ws.Range("N5").Value = 20.0924232795031
$ws.Range("O5").Value = 26.69919222243949

$ws.Range("B6").Value = 12.76078107049613
$ws.Range("C6").Value = 6.40571104506154
$ws.Range("D6").Value = 5.786406520101447
$ws.Range("E6").Value = 11.13630032507324
$ws.Range("G6").Value = 3.672325587855898
$ws.Range("I6").Value = 26.01281707464449
$ws.Range("K6").Value = 10.37124575189972
$ws.Range("L6").Value = 9.911632069085254
$ws.Range("N6").Value = 20.09519964894864
$ws.Range("O6").Value = 26.70225532282979

$ws.Range("B7").Value = 12.82783161070862
$ws.Range("C7").Value = 6.428094195162432
$ws.Range("D7").Value = 5.819143802182893
$ws.Range("E7").Value = 11.13254174608505
$ws.Range("G7").Value = 3.671705571697975
$ws.Range("I7").Value = 25.98794038671728
$ws.Range("K7").Value = 10.41734623469662
$ws.Range("L7").Value = 9.914033980851594
$ws.Range("N7").Value = 20.07608595902377
$ws.Range("O7").Value = 26.68146202812207

$ws.Range("B8").Value = 13.12543681128355
$ws.Range("C8").Value = 6.524321022792989
$ws.Range("D8").Value = 5.965213921975632
$ws.Range("E8").Value = 11.12015098823453
$ws.Range("G8").Value = 3.669113409647975
$ws.Range("I8").Value = 25.8882733368217
$ws.Range("K8").Value = 10.62280212708372
$ws.Range("L8").Value = 9.92851458012243
$ws.Range("N8").Value = 19.99562597494034
$ws.Range("O8").Value = 26.60130105863525

$ws.Range("B9").Value = 13.71003694709947
$ws.Range("C9").Value = 6.70332890909608
$ws.Range("D9").Value = 6.254575977444935
$ws.Range("E9").Value = 11.11089878540989
$ws.Range("G9").Value = 3.664543688964805
$ws.Range("I9").Value = 25.72924238850399
$ws.Range("K9").Value = 11.02966720059623
$ws.Range("L9").Value = 9.970700219672302
$ws.Range("N9").Value = 19.85171773084753
$ws.Range("O9").Value = 26.48591455874489

$ws.Range("B10").Value = 14.13565009265682
$ws.Range("C10").Value = 6.828474234138991
$ws.Range("D10").Value = 6.466474150904282
$ws.Range("E10").Value = 11.11323629695878
$ws.Range("G10").Value = 3.661496084322002
$ws.Range("I10").Value = 25.63467694484628
$ws.Range("K10").Value = 11.3280961976613
$ws.Range("L10").Value = 10.00993427940919
$ws.Range("N10").Value = 19.75436043594757
$ws.Range("O10").Value = 26.42673652656856

$ws.Range("B11").Value = 14.32752325042197
$ws.Range("C11").Value = 6.88392924358253
$ws.Range("D11").Value = 6.562224548133661
$ws.Range("E11").Value = 11.11627227605846
$ws.Range("G11").Value = 3.660176217233836
$ws.Range("I11").Value = 25.59650816719402
$ws.Range("K11").Value = 11.46315462872333
$ws.Range("L11").Value = 10.02953091184631
$ws.Range("N11").Value = 19.71186972341364
$ws.Range("O11").Value = 26.40539143712245

$ws.Range("B12").Value = 14.39985917934612
$ws.Range("C12").Value = 6.904709057930668
$ws.Range("D12").Value = 6.598351331725684
$ws.Range("E12").Value = 11.11770435877086
$ws.Range("G12").Value = 3.659685928279463
$ws.Range("I12").Value = 25.5827531199699
$ws.Range("K12").Value = 11.51414944092846
$ws.Range("L12").Value = 10.03719916211349
$ws.Range("N12").Value = 19.69603674434509
$ws.Range("O12").Value = 26.39811138109778

$ws.Range("B13").Value = 14.38429578433744
$ws.Range("C13").Value = 6.900243649703902
$ws.Range("D13").Value = 6.59057724904765
$ws.Range("E13").Value = 11.11738339475092
$ws.Range("G13").Value = 3.659791098293945
$ws.Range("I13").Value = 25.58568442368499
$ws.Range("K13").Value = 11.50317416938117
$ws.Range("L13").Value = 10.0355367301148
$ws.Range("N13").Value = 19.69943523317623
$ws.Range("O13").Value = 26.39964354845415

$ws.Range("B14").Value = 14.33348122937028
$ws.Range("C14").Value = 6.885643260091072
$ws.Range("D14").Value = 6.565199580003831
$ws.Range("E14").Value = 11.11638444185785
$ws.Range("G14").Value = 3.660135690418605
$ws.Range("I14").Value = 25.5953625223156
$ws.Range("K14").Value = 11.46735325752224
$ws.Range("L14").Value = 10.03015684871836
$ws.Range("N14").Value = 19.71056198400845
$ws.Range("O14").Value = 26.4047764063277

$ws.Range("B15").Value = 14.30231180172915
$ws.Range("C15").Value = 6.876671249500334
$ws.Range("D15").Value = 6.549636726981249
$ws.Range("E15").Value = 11.11580929466252
$ws.Range("G15").Value = 3.66034800092757
$ws.Range("I15").Value = 25.60138165753583
$ws.Range("K15").Value = 11.44539111855596
$ws.Range("L15").Value = 10.02689361384334
$ws.Range("N15").Value = 19.71741091802634
$ws.Range("O15").Value = 26.40802501610149

$ws.Range("B16").Value = 14.12306947358069
$ws.Range("C16").Value = 6.824819780211056
$ws.Range("D16").Value = 6.460200286661093
$ws.Range("E16").Value = 11.11307748083168
$ws.Range("G16").Value = 3.661583674929934
$ws.Range("I16").Value = 25.63726906071125
$ws.Range("K16").Value = 11.31925157487817
$ws.Range("L16").Value = 10.00868844477658
$ws.Range("N16").Value = 19.75717337681797
$ws.Range("O16").Value = 26.42824377178423

$ws.Range("B17").Value = 14.01261201123608
$ws.Range("C17").Value = 6.792627357538171
$ws.Range("D17").Value = 6.405140573693832
$ws.Range("E17").Value = 11.11190601056767
$ws.Range("G17").Value = 3.662358720453152
$ws.Range("I17").Value = 25.66052783875129
$ws.Range("K17").Value = 11.24165513454019
$ws.Range("L17").Value = 9.997965292742496
$ws.Range("N17").Value = 19.78202593771188
$ws.Range("O17").Value = 26.442076192821

$ws.Range("B18").Value = 13.94892007783745
$ws.Range("C18").Value = 6.773972815309784
$ws.Range("D18").Value = 6.373413375134257
$ws.Range("E18").Value = 11.11141796993092
$ws.Range("G18").Value = 3.662810768556477
$ws.Range("I18").Value = 25.67436206033011
$ws.Range("K18").Value = 11.19696076832985
$ws.Range("L18").Value = 9.991962490831535
$ws.Range("N18").Value = 19.7964897416627
$ws.Range("O18").Value = 26.45055695734265

$ws.Range("B19").Value = 13.92732990262699
$ws.Range("C19").Value = 6.767633194954342
$ws.Range("D19").Value = 6.36266230824397
$ws.Range("E19").Value = 11.11128466695794
$ws.Range("G19").Value = 3.662964901293813
$ws.Range("I19").Value = 25.67912443456599
$ws.Range("K19").Value = 11.18181871831727
$ws.Range("L19").Value = 9.989958481060352
$ws.Range("N19").Value = 19.80141604475592
$ws.Range("O19").Value = 26.45351848056228

$ws.Range("B20").Value = 14.02438746016731
$ws.Range("C20").Value = 6.796068668101737
$ws.Range("D20").Value = 6.411008100228148
$ws.Range("E20").Value = 11.11201149877518
$ws.Range("G20").Value = 3.662275567739737
$ws.Range("I20").Value = 25.65800465812222
$ws.Range("K20").Value = 11.24992227944855
$ws.Range("L20").Value = 9.999089755487777
$ws.Range("N20").Value = 19.77936282999112
$ws.Range("O20").Value = 26.44054939205245

$ws.Range("B21").Value = 14.34841599841732
$ws.Range("C21").Value = 6.889937774351501
$ws.Range("D21").Value = 6.572657497101318
$ws.Range("E21").Value = 11.11667020433386
$ws.Range("G21").Value = 3.660034217467854
$ws.Range("I21").Value = 25.59250085940918
$ws.Range("K21").Value = 11.47787914541955
$ws.Range("L21").Value = 10.03173036875379
$ws.Range("N21").Value = 19.70728681287668
$ws.Range("O21").Value = 26.4032469642423

$ws.Range("B22").Value = 14.55827809890095
$ws.Range("C22").Value = 6.950002117585585
$ws.Range("D22").Value = 6.677519659508836
$ws.Range("E22").Value = 11.1213604268444
$ws.Range("G22").Value = 3.658624810852758
$ws.Range("I22").Value = 25.55376308877066
$ws.Range("K22").Value = 11.62597317491933
$ws.Range("L22").Value = 10.05450304134397
$ws.Range("N22").Value = 19.66168036289707
$ws.Range("O22").Value = 26.38354772267386

$ws.Range("B23").Value = 14.44646830077507
$ws.Range("C23").Value = 6.918064686266044
$ws.Range("D23").Value = 6.621637034305121
$ws.Range("E23").Value = 11.11870705597561
$ws.Range("G23").Value = 3.659371979885043
$ws.Range("I23").Value = 25.57406510614664
$ws.Range("K23").Value = 11.54702947093711
$ws.Range("L23").Value = 10.04221847343934
$ws.Range("N23").Value = 19.68588457020549
$ws.Range("O23").Value = 26.39363304746697

$ws.Range("B24").Value = 14.01906436120884
$ws.Range("C24").Value = 6.794513307932207
$ws.Range("D24").Value = 6.408355614728412
$ws.Range("E24").Value = 11.11196322980583
$ws.Range("G24").Value = 3.662313140942534
$ws.Range("I24").Value = 25.65914394768768
$ws.Range("K24").Value = 11.24618495847542
$ws.Range("L24").Value = 9.998580880715657
$ws.Range("N24").Value = 19.78056627347723
$ws.Range("O24").Value = 26.4412380130572

$ws.Range("B25").Value = 13.55225399610933
$ws.Range("C25").Value = 6.655988887195205
$ws.Range("D25").Value = 6.176243990728836
$ws.Range("E25").Value = 11.11179421377089
$ws.Range("G25").Value = 3.665725285610473
$ws.Range("I25").Value = 25.76835794819988
$ws.Range("K25").Value = 11.02966720059623
$ws.Range("L25").Value = 9.9578276353755
$ws.Range("N25").Value = 19.88917235302576
$ws.Range("O25").Value = 26.51264156250693

